# update code tinh luong cho Quyen
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lương")

$ws.Range("B2").Value = 19
$ws.Range("B3").Value = 665000
$ws.Range("B4").Value = 2035714.285714285
$ws.Range("B12").Value = 1357142.857142857
$ws.Range("B20").Value = 2035714.285714285
$ws.Range("B28").Value = -839285.7142857146
$ws.Range("B29").Value = 1357142.857142857
$ws.Range("B30").Value = 2035714.285714285
$ws.Range("B31").Value = 2553571.428571428
